$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.053156097523423
$ws.Cells.Item(2, 4).Value = 1.051751861300523
$ws.Cells.Item(2, 5).Value = 1.066659961403967
$ws.Cells.Item(2, 6).Value = 1.073738682599607
$ws.Cells.Item(2, 9).Value = 1.046598637936877
$ws.Cells.Item(2, 10).Value = 1.058174923673622
$ws.Cells.Item(2, 11).Value = 1.054502493799493
$ws.Cells.Item(2, 12).Value = 1.069369913909654
$ws.Cells.Item(2, 13).Value = 1.076429746941927
$ws.Cells.Item(2, 14).Value = 1.023118874588622
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.054234686470764
$ws.Cells.Item(3, 4).Value = 1.052554283269969
$ws.Cells.Item(3, 5).Value = 1.067727661019617
$ws.Cells.Item(3, 6).Value = 1.074879733386872
$ws.Cells.Item(3, 9).Value = 1.046903821172418
$ws.Cells.Item(3, 10).Value = 1.058903655958421
$ws.Cells.Item(3, 11).Value = 1.055117651074186
$ws.Cells.Item(3, 12).Value = 1.070252620236714
$ws.Cells.Item(3, 13).Value = 1.077386993159354
$ws.Cells.Item(3, 14).Value = 1.023367733228577
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.054932768199299
$ws.Cells.Item(4, 4).Value = 1.053073590002316
$ws.Cells.Item(4, 5).Value = 1.068419062255575
$ws.Cells.Item(4, 6).Value = 1.075618732034853
$ws.Cells.Item(4, 9).Value = 1.047100207644079
$ws.Cells.Item(4, 10).Value = 1.059374765342897
$ws.Cells.Item(4, 11).Value = 1.055515130536403
$ws.Cells.Item(4, 12).Value = 1.07082371544046
$ws.Cells.Item(4, 13).Value = 1.078006467910669
$ws.Cells.Item(4, 14).Value = 1.023528469288137
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.055226280962694
$ws.Cells.Item(5, 4).Value = 1.053291926695964
$ws.Cells.Item(5, 5).Value = 1.068709853690404
$ws.Cells.Item(5, 6).Value = 1.075929565671266
$ws.Cells.Item(5, 9).Value = 1.047182508071778
$ws.Cells.Item(5, 10).Value = 1.059572716910298
$ws.Cells.Item(5, 11).Value = 1.055682094624096
$ws.Cells.Item(5, 12).Value = 1.071063785885057
$ws.Cells.Item(5, 13).Value = 1.078266912409261
$ws.Cells.Item(5, 14).Value = 1.023595972649722
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.055275565354677
$ws.Cells.Item(6, 4).Value = 1.053328587555422
$ws.Cells.Item(6, 5).Value = 1.068758686258729
$ws.Cells.Item(6, 6).Value = 1.075981765307427
$ws.Cells.Item(6, 9).Value = 1.047196311400312
$ws.Cells.Item(6, 10).Value = 1.05960594781787
$ws.Cells.Item(6, 11).Value = 1.055710120632048
$ws.Cells.Item(6, 12).Value = 1.071104093712859
$ws.Cells.Item(6, 13).Value = 1.078310643201082
$ws.Cells.Item(6, 14).Value = 1.023607302647599
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.054936689976783
$ws.Cells.Item(7, 4).Value = 1.053076507349326
$ws.Cells.Item(7, 5).Value = 1.068422947328576
$ws.Cells.Item(7, 6).Value = 1.07562288478589
$ws.Cells.Item(7, 9).Value = 1.047101308369909
$ws.Cells.Item(7, 10).Value = 1.059377410785528
$ws.Cells.Item(7, 11).Value = 1.055517362054434
$ws.Cells.Item(7, 12).Value = 1.070826923344212
$ws.Cells.Item(7, 13).Value = 1.078009947914803
$ws.Cells.Item(7, 14).Value = 1.023529371546894
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.053520578222193
$ws.Cells.Item(8, 4).Value = 1.052023025466011
$ws.Cells.Item(8, 5).Value = 1.067020686065584
$ws.Cells.Item(8, 6).Value = 1.074124168884341
$ws.Cells.Item(8, 9).Value = 1.046702001152556
$ws.Cells.Item(8, 10).Value = 1.058421290965924
$ws.Cells.Item(8, 11).Value = 1.05471050666004
$ws.Cells.Item(8, 12).Value = 1.069668244095064
$ws.Cells.Item(8, 13).Value = 1.076753237765805
$ws.Cells.Item(8, 14).Value = 1.023203037930606
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.051026440164228
$ws.Cells.Item(9, 4).Value = 1.050167334900684
$ws.Cells.Item(9, 5).Value = 1.064553771623398
$ws.Cells.Item(9, 6).Value = 1.07148831326722
$ws.Cells.Item(9, 9).Value = 1.045990046142234
$ws.Cells.Item(9, 10).Value = 1.056733202357038
$ws.Cells.Item(9, 11).Value = 1.05328437921718
$ws.Cells.Item(9, 12).Value = 1.067625932872791
$ws.Cells.Item(9, 13).Value = 1.074539309614892
$ws.Cells.Item(9, 14).Value = 1.022625765308756
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.04936448844024
$ws.Cells.Item(10, 4).Value = 1.048930689415009
$ws.Cells.Item(10, 5).Value = 1.062911886556176
$ws.Cells.Item(10, 6).Value = 1.069734483278302
$ws.Cells.Item(10, 9).Value = 1.045509816920129
$ws.Cells.Item(10, 10).Value = 1.055605604645871
$ws.Cells.Item(10, 11).Value = 1.052330715527157
$ws.Cells.Item(10, 12).Value = 1.066264009578558
$ws.Cells.Item(10, 13).Value = 1.073063729577532
$ws.Cells.Item(10, 14).Value = 1.02223942444673
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.048645030810528
$ws.Cells.Item(11, 4).Value = 1.048395327104036
$ws.Cells.Item(11, 5).Value = 1.06220157794824
$ws.Cells.Item(11, 6).Value = 1.068975861805396
$ws.Cells.Item(11, 9).Value = 1.045300546440236
$ws.Cells.Item(11, 10).Value = 1.055116819041291
$ws.Cells.Item(11, 11).Value = 1.051917079275075
$ws.Cells.Item(11, 12).Value = 1.065674190396395
$ws.Cells.Item(11, 13).Value = 1.072424874141129
$ws.Cells.Item(11, 14).Value = 1.022071781957639
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.04837781843299
$ws.Cells.Item(12, 4).Value = 1.048196487027724
$ws.Cells.Item(12, 5).Value = 1.061937833476979
$ws.Cells.Item(12, 6).Value = 1.068694195708288
$ws.Cells.Item(12, 9).Value = 1.045222614523619
$ws.Cells.Item(12, 10).Value = 1.054935182806051
$ws.Cells.Item(12, 11).Value = 1.051763332177892
$ws.Cells.Item(12, 12).Value = 1.065455090498199
$ws.Cells.Item(12, 13).Value = 1.072187586662736
$ws.Cells.Item(12, 14).Value = 1.022009458912581
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.048435135217931
$ws.Cells.Item(13, 4).Value = 1.04823913811045
$ws.Cells.Item(13, 5).Value = 1.061994403229846
$ws.Cells.Item(13, 6).Value = 1.068754608630052
$ws.Cells.Item(13, 9).Value = 1.045239340214253
$ws.Cells.Item(13, 10).Value = 1.054974147994457
$ws.Cells.Item(13, 11).Value = 1.051796316165886
$ws.Cells.Item(13, 12).Value = 1.065502088842002
$ws.Cells.Item(13, 13).Value = 1.072238485098317
$ws.Cells.Item(13, 14).Value = 1.022022829821701
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.048622942420545
$ws.Cells.Item(14, 4).Value = 1.048378890566213
$ws.Cells.Item(14, 5).Value = 1.062179774795098
$ws.Cells.Item(14, 6).Value = 1.068952576757725
$ws.Cells.Item(14, 9).Value = 1.045294108634395
$ws.Cells.Item(14, 10).Value = 1.055101806555849
$ws.Cells.Item(14, 11).Value = 1.051904372617058
$ws.Cells.Item(14, 12).Value = 1.065656079833119
$ws.Cells.Item(14, 13).Value = 1.072405259643095
$ws.Cells.Item(14, 14).Value = 1.022066631399306
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.048738660093399
$ws.Cells.Item(15, 4).Value = 1.048464998970887
$ws.Cells.Item(15, 5).Value = 1.062294001020814
$ws.Cells.Item(15, 6).Value = 1.069074567307828
$ws.Cells.Item(15, 9).Value = 1.045327826823833
$ws.Cells.Item(15, 10).Value = 1.055180450670938
$ws.Cells.Item(15, 11).Value = 1.051970935959434
$ws.Cells.Item(15, 12).Value = 1.065750956800307
$ws.Cells.Item(15, 13).Value = 1.07250801651175
$ws.Cells.Item(15, 14).Value = 1.022093611954925
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.049412240182784
$ws.Cells.Item(16, 4).Value = 1.04896622206158
$ws.Cells.Item(16, 5).Value = 1.062959040837975
$ws.Cells.Item(16, 6).Value = 1.069784847287236
$ws.Cells.Item(16, 9).Value = 1.045523677530513
$ws.Cells.Item(16, 10).Value = 1.05563803259462
$ws.Cells.Item(16, 11).Value = 1.05235815259222
$ws.Cells.Item(16, 12).Value = 1.066303151869351
$ws.Cells.Item(16, 13).Value = 1.073106129985797
$ws.Cells.Item(16, 14).Value = 1.022250542867448
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.049834806533041
$ws.Cells.Item(17, 4).Value = 1.049280656617816
$ws.Cells.Item(17, 5).Value = 1.063376373573363
$ws.Cells.Item(17, 6).Value = 1.070230600740441
$ws.Cells.Item(17, 9).Value = 1.045646173933162
$ws.Cells.Item(17, 10).Value = 1.055924920106034
$ws.Cells.Item(17, 11).Value = 1.052600857683073
$ws.Cells.Item(17, 12).Value = 1.066649503056461
$ws.Cells.Item(17, 13).Value = 1.073481332156404
$ws.Cells.Item(17, 14).Value = 1.022348886668107
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.050081299591404
$ws.Cells.Item(18, 4).Value = 1.049464071891109
$ws.Cells.Item(18, 5).Value = 1.063619858309347
$ws.Cells.Item(18, 6).Value = 1.070490678248278
$ws.Cells.Item(18, 9).Value = 1.045717495875064
$ws.Cells.Item(18, 10).Value = 1.05609220576703
$ws.Cells.Item(18, 11).Value = 1.052742356533955
$ws.Cells.Item(18, 12).Value = 1.066851514456072
$ws.Cells.Item(18, 13).Value = 1.073700189011635
$ws.Cells.Item(18, 14).Value = 1.022406214748802
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.050165350299121
$ws.Cells.Item(19, 4).Value = 1.049526613602086
$ws.Cells.Item(19, 5).Value = 1.063702890767166
$ws.Cells.Item(19, 6).Value = 1.070579371036369
$ws.Cells.Item(19, 9).Value = 1.045741793080369
$ws.Cells.Item(19, 10).Value = 1.056149237212525
$ws.Cells.Item(19, 11).Value = 1.052790592633457
$ws.Cells.Item(19, 12).Value = 1.066920393535352
$ws.Cells.Item(19, 13).Value = 1.073774814869109
$ws.Cells.Item(19, 14).Value = 1.022425756340546
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.049789467382071
$ws.Cells.Item(20, 4).Value = 1.049246919634365
$ws.Cells.Item(20, 5).Value = 1.063331591348581
$ws.Cells.Item(20, 6).Value = 1.070182767644596
$ws.Cells.Item(20, 9).Value = 1.045633044482369
$ws.Cells.Item(20, 10).Value = 1.055894145064038
$ws.Cells.Item(20, 11).Value = 1.052574824645586
$ws.Cells.Item(20, 12).Value = 1.066612343824
$ws.Cells.Item(20, 13).Value = 1.073441075711378
$ws.Cells.Item(20, 14).Value = 1.022338338845317
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.048567637161246
$ws.Cells.Item(21, 4).Value = 1.048337736470669
$ws.Cells.Item(21, 5).Value = 1.062125184846802
$ws.Cells.Item(21, 6).Value = 1.068894276776305
$ws.Cells.Item(21, 9).Value = 1.045277986208208
$ws.Cells.Item(21, 10).Value = 1.055064216486821
$ws.Cells.Item(21, 11).Value = 1.051872555560971
$ws.Cells.Item(21, 12).Value = 1.065610733733953
$ws.Cells.Item(21, 13).Value = 1.072356148370265
$ws.Cells.Item(21, 14).Value = 1.022053734393301
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.047799575420076
$ws.Cells.Item(22, 4).Value = 1.047766197301264
$ws.Cells.Item(22, 5).Value = 1.061367223233396
$ws.Cells.Item(22, 6).Value = 1.068084843818614
$ws.Cells.Item(22, 9).Value = 1.045053592480268
$ws.Cells.Item(22, 10).Value = 1.05454194743921
$ws.Cells.Item(22, 11).Value = 1.051430408398524
$ws.Cells.Item(22, 12).Value = 1.064980896056974
$ws.Cells.Item(22, 13).Value = 1.071674080096757
$ws.Cells.Item(22, 14).Value = 1.021874484746833
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.04820672532429
$ws.Cells.Item(23, 4).Value = 1.048069171416701
$ws.Cells.Item(23, 5).Value = 1.06176898055881
$ws.Cells.Item(23, 6).Value = 1.068513874002476
$ws.Cells.Item(23, 9).Value = 1.045172657279352
$ws.Cells.Item(23, 10).Value = 1.054818855786067
$ws.Cells.Item(23, 11).Value = 1.051664856096306
$ws.Cells.Item(23, 12).Value = 1.065314793098076
$ws.Cells.Item(23, 13).Value = 1.0720356509607
$ws.Cells.Item(23, 14).Value = 1.021969537500962
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.049809954139235
$ws.Cells.Item(24, 4).Value = 1.049262163891207
$ws.Cells.Item(24, 5).Value = 1.063351826317803
$ws.Cells.Item(24, 6).Value = 1.070204381123369
$ws.Cells.Item(24, 9).Value = 1.04563897751189
$ws.Cells.Item(24, 10).Value = 1.055908051138181
$ws.Cells.Item(24, 11).Value = 1.052586588061334
$ws.Cells.Item(24, 12).Value = 1.066629134510507
$ws.Cells.Item(24, 13).Value = 1.073459265841872
$ws.Cells.Item(24, 14).Value = 1.022343105058039
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.051671090860145
$ws.Cells.Item(25, 4).Value = 1.050646992339617
$ws.Cells.Item(25, 5).Value = 1.065191048083169
$ws.Cells.Item(25, 6).Value = 1.072169144325822
$ws.Cells.Item(25, 9).Value = 1.04617508963862
$ws.Cells.Item(25, 10).Value = 1.057170002593507
$ws.Cells.Item(25, 11).Value = 1.053653581305161
$ws.Cells.Item(25, 12).Value = 1.068153986771128
$ws.Cells.Item(25, 13).Value = 1.075111597424884
$ws.Cells.Item(25, 14).Value = 1.022775267729349
